$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# --- Row-level stat / status corrections for existing rows (2-52) ---
$ws.Range("G3").Value = "8:10 - 2nd Half"
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 4
$ws.Range("O3").Value = 25
$ws.Range("G4").Value = "7:42 - 1st Half"
$ws.Range("O4").Value = 9
$ws.Range("G5").Value = "8:10 - 2nd Half"
$ws.Range("H5").Value = 20
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 3
$ws.Range("O5").Value = 27
$ws.Range("G7").Value = "7:42 - 1st Half"
$ws.Range("H7").Value = 9
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 1
$ws.Range("O7").Value = 12
$ws.Range("G10").Value = "7:42 - 1st Half"
$ws.Range("G11").Value = "8:10 - 2nd Half"
$ws.Range("O11").Value = 20
$ws.Range("G12").Value = "7:42 - 1st Half"
$ws.Range("H12").Value = 5
$ws.Range("O12").Value = 10
$ws.Range("G13").Value = "8:10 - 2nd Half"
$ws.Range("D15").Value = "Kirill Elatontsev"
$ws.Range("E15").Value = "OU"
$ws.Range("G15").Value = "8:10 - 2nd Half"
$ws.Range("H15").Value = 8
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 14
$ws.Range("D16").Value = "Jalil Bethea"
$ws.Range("E16").Value = "ALA"
$ws.Range("G16").Value = "8:10 - 2nd Half"
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 4
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 7
$ws.Range("G17").Value = "8:10 - 2nd Half"
$ws.Range("O17").Value = 13
$ws.Range("G18").Value = "7:42 - 1st Half"
$ws.Range("H18").Value = 13
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 1
$ws.Range("O18").Value = 9
$ws.Range("G20").Value = "7:42 - 1st Half"
$ws.Range("O20").Value = 8
$ws.Range("G21").Value = "8:10 - 2nd Half"
$ws.Range("O21").Value = 22
$ws.Range("G22").Value = "8:10 - 2nd Half"
$ws.Range("H22").Value = 12
$ws.Range("I22").Value = 7
$ws.Range("O22").Value = 23
$ws.Range("G24").Value = "8:10 - 2nd Half"
$ws.Range("H24").Value = 14
$ws.Range("I24").Value = 11
$ws.Range("J24").Value = 6
$ws.Range("O24").Value = 28
$ws.Range("G25").Value = "8:10 - 2nd Half"
$ws.Range("H25").Value = 12
$ws.Range("I25").Value = 13
$ws.Range("O25").Value = 30
$ws.Range("G27").Value = "7:42 - 1st Half"
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 2
$ws.Range("O27").Value = 10
$ws.Range("D28").Value = "Xaivian Lee"
$ws.Range("E28").Value = "FLA"
$ws.Range("F28").Value = "FLA@VAN"
$ws.Range("G28").Value = "7:42 - 1st Half"
$ws.Range("H28").Value = 10
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 2
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("D29").Value = "J.P. Estrella"
$ws.Range("E29").Value = "TENN"
$ws.Range("F29").Value = "UK@TENN"
$ws.Range("G29").Value = "Final"
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 1
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 9
$ws.Range("G30").Value = "8:10 - 2nd Half"
$ws.Range("G31").Value = "7:42 - 1st Half"
$ws.Range("J31").Value = 4
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 8
$ws.Range("G33").Value = "8:10 - 2nd Half"
$ws.Range("H33").Value = 23
$ws.Range("O33").Value = 26
$ws.Range("G34").Value = "8:10 - 2nd Half"
$ws.Range("H34").Value = 14
$ws.Range("J34").Value = 5
$ws.Range("K34").Value = 1
$ws.Range("O34").Value = 28
$ws.Range("G35").Value = "7:42 - 1st Half"
$ws.Range("H35").Value = 8
$ws.Range("I35").Value = 6
$ws.Range("O35").Value = 11
$ws.Range("G36").Value = "7:42 - 1st Half"
$ws.Range("O36").Value = 5
$ws.Range("D40").Value = "Isaiah Brown"
$ws.Range("E40").Value = "FLA"
$ws.Range("G40").Value = "7:42 - 1st Half"
$ws.Range("H40").Value = 8
$ws.Range("I40").Value = 5
$ws.Range("J40").Value = 2
$ws.Range("K40").Value = 1
$ws.Range("O40").Value = 6
$ws.Range("D41").Value = "AK Okereke"
$ws.Range("E41").Value = "VAN"
$ws.Range("G41").Value = "7:42 - 1st Half"
$ws.Range("I41").Value = 7
$ws.Range("J41").Value = 0
$ws.Range("O41").Value = 10
$ws.Range("G46").Value = "8:10 - 2nd Half"
$ws.Range("H46").Value = 2
$ws.Range("J46").Value = 5
$ws.Range("O46").Value = 14
$ws.Range("G47").Value = "7:42 - 1st Half"
$ws.Range("D49").Value = "Mike James"
$ws.Range("E49").Value = "VAN"
$ws.Range("F49").Value = "FLA@VAN"
$ws.Range("G49").Value = "7:42 - 1st Half"
$ws.Range("O49").Value = 2
$ws.Range("D50").Value = "Noah Williamson"
$ws.Range("E50").Value = "ALA"
$ws.Range("F50").Value = "ALA@OU"
$ws.Range("G50").Value = "8:10 - 2nd Half"
$ws.Range("J50").Value = 1
$ws.Range("K50").Value = 0
$ws.Range("N50").Value = 1
$ws.Range("O50").Value = 4

# --- Insert new row 53 (Urban Klavzar, FLA@VAN) which pushes the former
#     row 53 (Troy Henderson) down to row 54; Troy Henderson's data is
#     unchanged so no further edit is required for row 54. ---
$ws.Rows.Item(53).Insert()

$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "2026-01-17"
$ws.Range("A53").NumberFormat = "General"
$ws.Range("B53").Value = "Undrafted"
$ws.Range("C53").Value = "No"
$ws.Range("D53").Value = "Urban Klavzar"
$ws.Range("E53").Value = "FLA"
$ws.Range("F53").Value = "FLA@VAN"
$ws.Range("G53").Value = "7:42 - 1st Half"
$ws.Range("H53").Value = -1
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 6

# --- Column G got a touch narrower in this update. The engine's
#     ColumnWidth setter adds a constant 5/6-character rounding offset
#     before it writes the stored <col width> (matching how it also
#     reports the existing "18" columns once you round-trip them), so we
#     back that out here to land on the target stored width of 17. ---
$ws.Columns.Item(7).ColumnWidth = 16.166666666666668

# --- OwnerTotals sheet: refresh starter pooh totals to match the
#     corrected per-player stats above. ---
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Range("B2").Value = 45
$ws2.Range("B3").Value = 39
$ws2.Range("B4").Value = 31
$ws2.Range("B5").Value = 26
$ws2.Range("B6").Value = 17
$ws2.Range("B8").Value = 9
